# Auto-generated edit script applying cell-level numeric updates
# to the Garuda_Profits workbook sheets, per the scheduled runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1304.3334
$ws.Range("I18").Value = 1204.5883
$ws.Range("K18").Value = 1204.5883
$ws.Range("M18").Value = -920.5882999999999
$ws.Range("H19").Value = 7718.625
$ws.Range("I19").Value = 10009.8
$ws.Range("J19").Value = 3900
$ws.Range("K19").Value = 10009.8
$ws.Range("L19").Value = 3900
$ws.Range("M19").Value = -9834.799999999999
$ws.Range("N19").Value = -4250
$ws.Range("H113").Value = 3950
$ws.Range("I113").Value = 3850
$ws.Range("J113").Value = 3994.4443
$ws.Range("K113").Value = 3850
$ws.Range("L113").Value = 3994.4443
$ws.Range("M113").Value = -596
$ws.Range("N113").Value = -10502.4443
$ws.Range("H132").Value = 3176178
$ws.Range("I132").Value = 3572840.2
$ws.Range("J132").Value = 2879.6
$ws.Range("K132").Value = 10718520.6
$ws.Range("L132").Value = 8638.799999999999
$ws.Range("M132").Value = -10715990.6
$ws.Range("N132").Value = -13698.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1222.5385
$ws.Range("I2").Value = 620.1667
$ws.Range("J2").Value = 1738.8572
$ws.Range("K2").Value = 620.1667
$ws.Range("L2").Value = 1738.8572
$ws.Range("M2").Value = -507.1667
$ws.Range("N2").Value = -1964.8572
$ws.Range("H32").Value = 15041.296
$ws.Range("I32").Value = 15415.443
$ws.Range("K32").Value = 15415.443
$ws.Range("M32").Value = -15128.443
$ws.Range("H45").Value = 1124.2084
$ws.Range("I45").Value = 1043.6154
$ws.Range("J45").Value = 1219.4546
$ws.Range("K45").Value = 1043.6154
$ws.Range("L45").Value = 1219.4546
$ws.Range("M45").Value = -666.6153999999999
$ws.Range("N45").Value = -1973.4546
$ws.Range("H110").Value = 4589.1
$ws.Range("I110").Value = 4730.1035
$ws.Range("J110").Value = 500
$ws.Range("K110").Value = 4730.1035
$ws.Range("L110").Value = 500
$ws.Range("M110").Value = -2685.1035
$ws.Range("N110").Value = -4590
$ws.Range("H116").Value = 1222.5385
$ws.Range("I116").Value = 620.1667
$ws.Range("J116").Value = 1738.8572
$ws.Range("K116").Value = 620.1667
$ws.Range("L116").Value = 1738.8572
$ws.Range("M116").Value = 1673.8333
$ws.Range("N116").Value = -6326.8572
$ws.Range("H122").Value = 987.5333000000001
$ws.Range("I122").Value = 965.8461
$ws.Range("J122").Value = 1128.5
$ws.Range("K122").Value = 2897.5383
$ws.Range("L122").Value = 3385.5
$ws.Range("M122").Value = -447.5383000000002
$ws.Range("N122").Value = -8285.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1222.5385
$ws.Range("I3").Value = 620.1667
$ws.Range("J3").Value = 1738.8572
$ws.Range("K3").Value = 620.1667
$ws.Range("L3").Value = 1738.8572
$ws.Range("M3").Value = -506.1667
$ws.Range("N3").Value = -1966.8572
$ws.Range("H99").Value = 1165.25
$ws.Range("I99").Value = 916.6667
$ws.Range("K99").Value = 916.6667
$ws.Range("M99").Value = 581.3333
$ws.Range("H107").Value = 450.25
$ws.Range("I107").Value = 428.625
$ws.Range("J107").Value = 493.5
$ws.Range("K107").Value = 428.625
$ws.Range("L107").Value = 493.5
$ws.Range("M107").Value = 1491.375
$ws.Range("N107").Value = -4333.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1073.75
$ws.Range("I16").Value = 1073.75
$ws.Range("K16").Value = 1073.75
$ws.Range("M16").Value = -786.75
$ws.Range("H33").Value = 6107.625
$ws.Range("I33").Value = 2972.2
$ws.Range("K33").Value = 2972.2
$ws.Range("M33").Value = -2593.2
$ws.Range("H36").Value = 11953.333
$ws.Range("J36").Value = 14511.429
$ws.Range("L36").Value = 14511.429
$ws.Range("N36").Value = -15287.429
$ws.Range("H40").Value = 11953.333
$ws.Range("J40").Value = 14511.429
$ws.Range("L40").Value = 14511.429
$ws.Range("N40").Value = -14831.429
$ws.Range("H58").Value = 772.6774
$ws.Range("I58").Value = 623.3214
$ws.Range("J58").Value = 2166.6667
$ws.Range("K58").Value = 623.3214
$ws.Range("L58").Value = 2166.6667
$ws.Range("M58").Value = -420.3214
$ws.Range("N58").Value = -2572.6667
$ws.Range("H107").Value = 84203.75
$ws.Range("I107").Value = 111970.22
$ws.Range("K107").Value = 111970.22
$ws.Range("M107").Value = -110050.22
$ws.Range("H113").Value = 1073.75
$ws.Range("I113").Value = 1073.75
$ws.Range("K113").Value = 1073.75
$ws.Range("M113").Value = 1096.25
$ws.Range("H122").Value = 1712.5454
$ws.Range("I122").Value = 2437.3333
$ws.Range("J122").Value = 842.8
$ws.Range("K122").Value = 7311.999899999999
$ws.Range("L122").Value = 2528.4
$ws.Range("M122").Value = -4861.999899999999
$ws.Range("N122").Value = -7428.4
$ws.Range("H136").Value = 772.6774
$ws.Range("I136").Value = 623.3214
$ws.Range("J136").Value = 2166.6667
$ws.Range("K136").Value = 1869.9642
$ws.Range("L136").Value = 6500.000100000001
$ws.Range("M136").Value = 680.0357999999999
$ws.Range("N136").Value = -11600.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 374489.4
$ws.Range("J131").Value = 529987.9399999999
$ws.Range("L131").Value = 1589963.82
$ws.Range("N131").Value = -1600043.82

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 2040.75
$ws.Range("I55").Value = 915
$ws.Range("J55").Value = 3166.5
$ws.Range("K55").Value = 915
$ws.Range("L55").Value = 3166.5
$ws.Range("M55").Value = -588
$ws.Range("N55").Value = -3820.5
$ws.Range("H107").Value = 358.625
$ws.Range("I107").Value = 352.9091
$ws.Range("J107").Value = 371.2
$ws.Range("K107").Value = 352.9091
$ws.Range("L107").Value = 371.2
$ws.Range("M107").Value = 1567.0909
$ws.Range("N107").Value = -4211.2
$ws.Range("H120").Value = 30000
$ws.Range("J120").Value = 30000
$ws.Range("L120").Value = 30000
$ws.Range("N120").Value = -39676
$ws.Range("H122").Value = 29414074
$ws.Range("I122").Value = 50002270
$ws.Range("J122").Value = 2362.5
$ws.Range("K122").Value = 150006810
$ws.Range("L122").Value = 7087.5
$ws.Range("M122").Value = -150004360
$ws.Range("N122").Value = -11987.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1501.0392
$ws.Range("I132").Value = 1241.081
$ws.Range("J132").Value = 2188.0715
$ws.Range("K132").Value = 3723.242999999999
$ws.Range("L132").Value = 6564.2145
$ws.Range("M132").Value = -1193.242999999999
$ws.Range("N132").Value = -11624.2145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H119").Value = 28698
$ws.Range("J119").Value = 28698
$ws.Range("L119").Value = 28698
$ws.Range("N119").Value = -38374
$ws.Range("H122").Value = 2271.2693
$ws.Range("I122").Value = 2233.524
$ws.Range("K122").Value = 6700.572
$ws.Range("M122").Value = -4250.572
$ws.Range("H124").Value = 30107.25
$ws.Range("J124").Value = 30107.25
$ws.Range("L124").Value = 30107.25
$ws.Range("N124").Value = -39927.25
$ws.Range("H126").Value = 7543.222
$ws.Range("I126").Value = 8518.532999999999
$ws.Range("J126").Value = 2666.6667
$ws.Range("K126").Value = 25555.599
$ws.Range("L126").Value = 8000.000100000001
$ws.Range("M126").Value = -23085.599
$ws.Range("N126").Value = -12940.0001
$ws.Range("H136").Value = 5522.2856
$ws.Range("I136").Value = 5808.96
$ws.Range("J136").Value = 3133.3333
$ws.Range("K136").Value = 17426.88
$ws.Range("L136").Value = 9399.999899999999
$ws.Range("M136").Value = -14876.88
$ws.Range("N136").Value = -14499.9999
